$d = $word.ActiveDocument

# Locate the old phone number text (the run containing " 07825 750236")
$rng = $d.Content
$rng.Find.Execute("07825 750236")

# Replace just the digits (leave the leading space's run untouched)
$rng.Text = "07875 725020"

# Force a run split at this exact boundary: briefly toggle a formatting
# property on the new text and back off again. Word (and this host) merge
# consecutive runs that share identical formatting, so a momentary Bold
# on/off here makes the corrected number live in its own <w:r>, matching
# how Word itself would split the run after an in-place retype.
$rng2 = $d.Content
$rng2.Find.Execute("07875 725020")
$rng2.Font.Bold = 1
$rng2.Font.Bold = 0
